$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.581.43"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "1.669.30"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'237.75"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4776"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("D8").Value = "'0.2617"
$ws.Range("E8").Value = "  +3.00%  "
$ws.Range("D9").Value = "'0.06172"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "1.669.91"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").Value = "'0.06992"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "'0.5895"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "'4.375"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "'75.39"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D16").Value = "'1.0000"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "25.572.03"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").Value = "'0.000006744"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("D20").Value = "'11.42"
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").Value = "1.886.18"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'4.445"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").Value = "'8.796"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'5.261"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "'136.82"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").Value = "'15.02"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").Value = "'1.383"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").Value = "'1.718"
$ws.Range("E28").Value = "  +5.32%  "
$ws.Range("D29").Value = "'104.75"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").Value = "'3.999"
$ws.Range("E30").Value = "  +6.83%  "
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("D32").Value = "'3.632"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04314"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.624"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'0.9559"
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6059"
$ws.Range("E36").Value = "  +5.13%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.583"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.9185"
$ws.Range("E38").Value = "  +12.20%  "
$ws.Range("B39").Value = "PaxDollar"
$ws.Range("C39").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D39").Value = "'0.9997"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.850"
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01472"
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'98.19"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.3762"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'4.884"
$ws.Range("E44").Value = "  +4.09%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.1121"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'6.217"
$ws.Range("E46").Value = "  +3.44%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05268"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'29.99"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.452"
$ws.Range("E49").Value = "  +4.85%  "
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").Value = "'1.002"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3344"
$ws.Range("E51").Value = "  +2.25%  "
